# Updated cryptos list with GitHub Actions: refresh Price (column D) and
# Volume(1h) (column E) figures for each coin row on the active sheet.
# Leading apostrophes force Excel to keep numeric-looking Price strings
# (e.g. "492.78") stored as text, matching the original cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "69.368.10"
$ws.Range('E2').Value = "  +1.68%  "
$ws.Range('D3').Value = "3.943.14"
$ws.Range('E3').Value = "  +0.53%  "
$ws.Range('E4').Value = "  +0.13%  "
$ws.Range('D5').Value = "'492.78"
$ws.Range('E5').Value = "  +0.86%  "
$ws.Range('D6').Value = "'147.73"
$ws.Range('E6').Value = "  +0.22%  "
$ws.Range('D7').Value = "'0.623"
$ws.Range('E7').Value = "  -0.75%  "
$ws.Range('E8').Value = "  -0.03%  "
$ws.Range('D9').Value = "'0.734"
$ws.Range('E9').Value = "  +0.30%  "
$ws.Range('D10').Value = "'0.177"
$ws.Range('E10').Value = "  +3.37%  "
$ws.Range('D11').Value = "'0.0000351"
$ws.Range('E11').Value = "  -1.42%  "
$ws.Range('D12').Value = "'43.29"
$ws.Range('E12').Value = "  +1.27%  "
$ws.Range('D13').Value = "'10.44"
$ws.Range('E13').Value = "  -1.71%  "
$ws.Range('D14').Value = "4.575.92"
$ws.Range('E14').Value = "  +0.67%  "
$ws.Range('D15').Value = "3.941.44"
$ws.Range('E15').Value = "  +0.50%  "
$ws.Range('D16').Value = "'14.32"
$ws.Range('E16').Value = "  -3.19%  "
$ws.Range('E18').Value = "  -0.70%  "
$ws.Range('E19').Value = "  +2.90%  "
$ws.Range('D20').Value = "69.439.17"
$ws.Range('E20').Value = "  +1.61%  "
$ws.Range('D21').Value = "'441.27"
$ws.Range('E21').Value = "  -0.63%  "
$ws.Range('E22').Value = "  +2.29%  "
$ws.Range('D23').Value = "'14.49"
$ws.Range('D24').Value = "'88.92"
$ws.Range('E24').Value = "  +0.46%  "
$ws.Range('D25').Value = "'12.10"
$ws.Range('E25').Value = "  +7.29%  "
$ws.Range('E26').Value = "  +4.28%  "
$ws.Range('D27').Value = "'11.11"
$ws.Range('E27').Value = "  -4.80%  "
$ws.Range('D28').Value = "'37.15"
$ws.Range('E28').Value = "  -4.49%  "
$ws.Range('E29').Value = "  -4.37%  "
$ws.Range('D30').Value = "'706.00"
$ws.Range('E30').Value = "  -1.16%  "
$ws.Range('E31').Value = "  -0.16%  "
$ws.Range('D32').Value = "'13.36"
$ws.Range('E32').Value = "  -0.83%  "
$ws.Range('D33').Value = "'2.89"
$ws.Range('D34').Value = "'0.465"
$ws.Range('E34').Value = "  +19.50%  "
$ws.Range('D35').Value = "0.0₃0914"
$ws.Range('E35').Value = "  +0.45%  "
$ws.Range('E36').Value = "  +3.42%  "
$ws.Range('D37').Value = "'61.56"
$ws.Range('E37').Value = "  +4.25%  "
$ws.Range('D38').Value = "'40.99"
$ws.Range('E38').Value = "  -0.84%  "
$ws.Range('E39').Value = "  +1.03%  "
$ws.Range('D40').Value = "'0.999"
$ws.Range('E40').Value = "  -0.12%  "
$ws.Range('E41').Value = "  +0.08%  "
$ws.Range('D42').Value = "'0.0490"
$ws.Range('E42').Value = "  +2.19%  "
$ws.Range('E43').Value = "  +0.42%  "
$ws.Range('E44').Value = "  -2.29%  "
$ws.Range('E45').Value = "  +2.81%  "
$ws.Range('D46').Value = "'0.143"
$ws.Range('E46').Value = "  +0.36%  "
$ws.Range('E47').Value = "  +7.28%  "
$ws.Range('D48').Value = "0.0₆0359"
$ws.Range('E48').Value = "  +4.84%  "
$ws.Range('D49').Value = "'3.01"
$ws.Range('E49').Value = "  +5.97%  "
$ws.Range('D50').Value = "'3.40"
$ws.Range('E50').Value = "  -0.57%  "
$ws.Range('D51').Value = "'144.14"
$ws.Range('E51').Value = "  -0.70%  "
